$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Scoot Henderson" (Portland Trail Blazers) row is removed, and
# several other rows shift position accordingly. Rewrite rows 2-18 with
# the final data, then clear the now-unused row 19.

$data = @(
    @("Collin Sexton", "PG,SG", "Utah Jazz"),
    @("James Harden", "PG,SG", "LA Clippers"),
    @("Fred VanVleet", "PG", "Houston Rockets"),
    @("D'Angelo Russell", "PG", "Los Angeles Lakers"),
    @("Anthony Edwards", "SG,SF", "Minnesota Timberwolves"),
    @("Jayson Tatum", "SF,PF", "Boston Celtics"),
    @("Jaren Jackson Jr.", "PF,C", "Memphis Grizzlies"),
    @("Kyle Kuzma", "PF", "Washington Wizards"),
    @("Ivica Zubac", "C", "LA Clippers"),
    @("Giannis Antetokounmpo", "PF,C", "Milwaukee Bucks"),
    @("Jrue Holiday", "PG,SG", "Boston Celtics"),
    @("Marcus Smart", "PG,SG", "Memphis Grizzlies"),
    @("Tyus Jones", "PG", "Phoenix Suns"),
    @("Jaden Ivey", "PG,SG", "Detroit Pistons"),
    @("Terry Rozier", "PG", "Miami Heat"),
    @("Wendell Carter Jr.", "C", "Orlando Magic"),
    @("Paul George", "SG,SF,PF", "Philadelphia 76ers")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}

# Remove the now-obsolete last row (previously row 19) entirely so the
# used range shrinks back to A1:C18.
$ws.Rows.Item(19).Delete()
